# Generate Report for Handback
# Adds a new handed-back file (b7047860-56d2-4699-9535-b3d37ad7c441.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" sheets/tables.

$wb = $excel.ActiveWorkbook

$fileName   = "b7047860-56d2-4699-9535-b3d37ad7c441.md"
$pathName   = "e2e\b7047860-56d2-4699-9535-b3d37ad7c441.md"
$extension  = ".md"
$status     = "Handed back: in sync with en-US"
$zhXlf      = "b7047860-56d2-4699-9535-b3d37ad7c441.9d27e32e940d9831fcc6e9dd213cdabf2462af55.zh-cn.xlf"
$deXlf      = "b7047860-56d2-4699-9535-b3d37ad7c441.9d27e32e940d9831fcc6e9dd213cdabf2462af55.de-de.xlf"
$dtFmt      = "yyyy-mm-dd HH:mm:ss"

# ===================== Overview =====================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$null = $loOv.ListRows.Add()

$wsOv.Range("A4").Value2 = $fileName
$wsOv.Range("C4").Value2 = $extension
$wsOv.Range("E4").Value2 = $status
$wsOv.Range("F4").Value2 = $status
$wsOv.Range("G4").NumberFormat = $dtFmt
$wsOv.Range("G4").Value2 = "2016-08-28 06:49:42"

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/changeset/e2e/b7047860-56d2-4699-9535-b3d37ad7c441.md", "", "", $pathName)

# ===================== zh-cn =====================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$null = $loZh.ListRows.Add()

$wsZh.Range("B4").Value2 = $extension
$wsZh.Range("C4").Value2 = $status
$wsZh.Range("D4").Value2 = "e2e"
$wsZh.Range("E4").Value2 = "ht"
$wsZh.Range("F4").Value2 = "'True"
$wsZh.Range("G4").Value2 = $zhXlf
$wsZh.Range("H4").NumberFormat = $dtFmt
$wsZh.Range("H4").Value2 = "2016-08-28 06:49:38"
$wsZh.Range("J4").Value2 = $zhXlf
$wsZh.Range("K4").NumberFormat = $dtFmt
$wsZh.Range("K4").Value2 = "2016-08-28 06:49:55"
$wsZh.Range("M4").Value2 = "'True"
$wsZh.Range("O4").Value2 = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/changeset/e2e/b7047860-56d2-4699-9535-b3d37ad7c441.md", "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/changeset/e2e/b7047860-56d2-4699-9535-b3d37ad7c441.md", "", "", $fileName)

# ===================== de-de =====================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$null = $loDe.ListRows.Add()

$wsDe.Range("B4").Value2 = $extension
$wsDe.Range("C4").Value2 = $status
$wsDe.Range("D4").Value2 = "e2e"
$wsDe.Range("E4").Value2 = "ht"
$wsDe.Range("F4").Value2 = "'True"
$wsDe.Range("G4").Value2 = $deXlf
$wsDe.Range("H4").NumberFormat = $dtFmt
$wsDe.Range("H4").Value2 = "2016-08-28 06:49:42"
$wsDe.Range("J4").Value2 = $deXlf
$wsDe.Range("K4").NumberFormat = $dtFmt
$wsDe.Range("K4").Value2 = "2016-08-28 06:50:05"
$wsDe.Range("M4").Value2 = "'True"
$wsDe.Range("O4").Value2 = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/changeset/e2e/b7047860-56d2-4699-9535-b3d37ad7c441.md", "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/changeset/e2e/b7047860-56d2-4699-9535-b3d37ad7c441.md", "", "", $fileName)

Write-Output "Handback report row added to Overview, zh-cn, de-de"
